# Commit: "Remove components in CONN1 and BAT"
#
# The BOM sheet has one row per placed component. Remove the BAT
# (battery header) row and the CONN1 (10-pin header) row entirely,
# exactly like selecting those sheet rows and clearing their contents
# in Excel: the row numbers of the surrounding rows are left
# untouched (row 8 / row 23 simply disappear from the saved XML
# because they end up with no cell data), matching how Excel omits
# fully-empty rows when it serialises the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CONN1 is on row 23 (Designator "CONN1", Device "HDR-M_2.54_1x10P").
$ws.Rows(23).Select()
$ws.Rows(23).ClearContents()

# BAT is on row 8 (Designator "BAT", Device "HDR-M-2.54_1x2").
# Selected/cleared last so it ends up as the sheet's final selection,
# matching the saved worksheet's <selection activeCell="A8" .../>.
$ws.Rows(8).Select()
$ws.Rows(8).ClearContents()
